$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")
$ws.Range("E10").Value = 6
$ws.Range("B11").Value = 21
